$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells -------------------------------------------------
# "<field>_old" -> "<field>_FV2304" and "<field>_new" -> "<field>_FV2310"
$ws.Range("A1").Value = "Segmentname_FV2304"
$ws.Range("B1").Value = "Segmentgruppe_FV2304"
$ws.Range("C1").Value = "Segment_FV2304"
$ws.Range("D1").Value = "Datenelement_FV2304"
$ws.Range("E1").Value = "Segment ID_FV2304"
$ws.Range("F1").Value = "Code_FV2304"
$ws.Range("G1").Value = "Qualifier_FV2304"
$ws.Range("H1").Value = "Beschreibung_FV2304"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("J1").Value = "Bedingung_FV2304"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2310"
$ws.Range("M1").Value = "Segmentgruppe_FV2310"
$ws.Range("N1").Value = "Segment_FV2310"
$ws.Range("O1").Value = "Datenelement_FV2310"
$ws.Range("P1").Value = "Segment ID_FV2310"
$ws.Range("Q1").Value = "Code_FV2310"
$ws.Range("R1").Value = "Qualifier_FV2310"
$ws.Range("S1").Value = "Beschreibung_FV2310"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("U1").Value = "Bedingung_FV2310"

# --- 2. Turn the data range into an Excel Table ("Table1") -----------------
# Preserve the header row's existing formatting across the table creation:
# stash it in an unused scratch row, clear the header formatting (so the new
# table does not derive a header-row dxf from it), build the table, then
# restore the original formatting onto the header cells.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A70:U70")
$headerRange.Copy()
$scratch.PasteSpecial(-4122)
$headerRange.ClearFormats()

$listObj = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), [System.Reflection.Missing]::Value, 1)
$listObj.Name = "Table1"

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$scratch.ClearContents()
$scratch.ClearFormats()

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
